# Weekly data refresh: insert one new "Choclo" (corn) price record.
#
# The new record is inserted as row 568, which pushes every existing
# record that used to live in rows 568-658 down by one row (to 569-659).
# Excel's native "insert entire row" semantics (shift-down + inherit the
# formatting of the surrounding rows) give us exactly that, so we let
# Excel do the shifting instead of rewriting every downstream row by hand.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 568..658 down to 569..659 and open up a blank row 568.
$ws.Rows("568").Insert()

# Populate the newly opened row 568 with the new record.
$ws.Range("A568").Value = 6
$ws.Range("B568").Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range("C568").Value = 'Metropolitana'
$ws.Range("D568").Value = 44505
$ws.Range("E568").Value = 13
$ws.Range("F568").Value = 100112024
$ws.Range("G568").Value = 'Choclo'
$ws.Range("H568").Value = 'Dulce o Americano'
$ws.Range("I568").Value = 'Primera'
$ws.Range("J568").Value = 350
$ws.Range("K568").Value = 18000
$ws.Range("L568").Value = 20000
$ws.Range("M568").Value = 19143
$ws.Range("N568").Value = '$/caja 50 unidades'
$ws.Range("O568").Value = 'Argentina'
$ws.Range("P568").Value = 383
$ws.Range("Q568").Value = 50
$ws.Range("R568").Value = 'Hortaliza'
